# Applies the 2016 campaign-files update: new party/candidate rows appended
# to the bottom of the sheet, plus translations filled in for an existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 59: fill in the previously-empty D/E translation cells ---------
$ws.Range("D59").Value = "Democratic Movement-United Georgia"
$ws.Range("E59").Value = "Democratic Movement-United Georgia"

# --- Existing rows 30, 41, 55 gain a numeric A-column id -----------------
$ws.Range("A30").Value = 36
$ws.Range("A41").Value = 38
$ws.Range("A55").Value = 37
$ws.Range("A59").Value = 34

# --- New rows 60-76: additional parties / candidates / blocs --------------
$ws.Range("A60").Value = 33
$ws.Range("C60").Value = "ირაკლი შიხიაშვილი დამოუკიდებელი კანდიდატი"
$ws.Range("D60").Value = "Irakli Shikhiashvili Independent Candidate"
$ws.Range("E60").Value = "Irakli Shikhiashvili Independent Candidate"

$ws.Range("A61").Value = 35
$ws.Range("C61").Value = "ზვიადის გზა – უფლის სახელით"
$ws.Range("D61").Value = "Zviadi’s Way – In the Name of the Lord"
$ws.Range("E61").Value = "Zviadi’s Way – In the Name of the Lord"

$ws.Range("A62").Value = 39
$ws.Range("C62").Value = "მოძრაობა-სახელმწიფო ხალხისთვის"
$ws.Range("D62").Value = "The movement- State for People"
$ws.Range("E62").Value = "The movement- State for People"

$ws.Range("A63").Value = 40
$ws.Range("C63").Value = "ქალთა პარტია"
$ws.Range("D63").Value = "Women's Party"
$ws.Range("E63").Value = "Women's Party"

$ws.Range("A64").Value = 41
$ws.Range("C64").Value = "ქართული იდეა"
$ws.Range("D64").Value = "Georgian Idea"
$ws.Range("E64").Value = "Georgian Idea"

$ws.Range("A65").Value = 42
$ws.Range("C65").Value = "წარმატებული საქართველო"
$ws.Range("D65").Value = "Successful Georgia"
$ws.Range("E65").Value = "Successful Georgia"

$ws.Range("A66").Value = 43
$ws.Range("C66").Value = "ჩვენი სამშობლო"
$ws.Range("D66").Value = "Our homeland"
$ws.Range("E66").Value = "Our homeland"

$ws.Range("A67").Value = 44
$ws.Range("C67").Value = 'მოქალაქეთა პოლიტიკური გაერთიანება "ახალი პოლიტიკური ცენტრი"'
$ws.Range("D67").Value = "New Political Center"
$ws.Range("E67").Value = "New Political Center"

$ws.Range("A68").Value = 45
$ws.Range("C68").Value = "სამოქალაქო პლატფორმა – ახალი საქართველო"
$ws.Range("D68").Value = "Civic Platform - New Georgia"
$ws.Range("E68").Value = "Civic Platform - New Georgia"

$ws.Range("A69").Value = 46
$ws.Range("C69").Value = "სალომე ზურაბიშვილი"
$ws.Range("D69").Value = "სალომე ზურაბიშვილი"
$ws.Range("E69").Value = "სალომე ზურაბიშვილი"

$ws.Range("A70").Value = 47
$ws.Range("C70").Value = 'მოქალაქეთა პოლიტიკური გაერთიანება  „ თოფაძე -მრეწველები“'
$ws.Range("D70").Value = 'მოქალაქეთა პოლიტიკური გაერთიანება  „ თოფაძე -მრეწველები“'
$ws.Range("E70").Value = 'მოქალაქეთა პოლიტიკური გაერთიანება  „ თოფაძე -მრეწველები“'

$ws.Range("A71").Value = 48
$ws.Range("C71").Value = "პლატფორმა ფინანსური სახელმწიფო ხალხისთვის"
$ws.Range("D71").Value = "პლატფორმა ფინანსური სახელმწიფო ხალხისთვის"
$ws.Range("E71").Value = "პლატფორმა ფინანსური სახელმწიფო ხალხისთვის"

$ws.Range("A72").Value = 49
$ws.Range("C72").Value = 'საარჩევნო ბლოკი   „ თოფაძე -მრეწველები, ჩვენი სამშობლო“'
$ws.Range("D72").Value = 'საარჩევნო ბლოკი   „ თოფაძე -მრეწველები, ჩვენი სამშობლო“'
$ws.Range("E72").Value = 'საარჩევნო ბლოკი   „ თოფაძე -მრეწველები, ჩვენი სამშობლო“'

$ws.Range("A73").Value = 50
$ws.Range("C73").Value = "მამული, ენა, სარწმუნოება"
$ws.Range("D73").Value = "მამული, ენა, სარწმუნოება"
$ws.Range("E73").Value = "მამული, ენა, სარწმუნოება"

$ws.Range("A74").Value = 51
$ws.Range("C74").Value = "ახალი პოლიტიკური მოძრაობა სახელმწიფო ხალხისთვის"
$ws.Range("D74").Value = "ახალი პოლიტიკური მოძრაობა სახელმწიფო ხალხისთვის"
$ws.Range("E74").Value = "ახალი პოლიტიკური მოძრაობა სახელმწიფო ხალხისთვის"

$ws.Range("A75").Value = 52
$ws.Range("C75").Value = "პროგრესულ-დემოკრატიული მოძრაობა"
$ws.Range("D75").Value = "პროგრესულ-დემოკრატიული მოძრაობა"
$ws.Range("E75").Value = "პროგრესულ-დემოკრატიული მოძრაობა"

$ws.Range("A76").Value = 53
$ws.Range("C76").Value = "საარჩევნო ბლოკი პაატა ბურჭულაძე სახელმწიფო ხალხისთვის"
$ws.Range("D76").Value = "საარჩევნო ბლოკი პაატა ბურჭულაძე სახელმწიფო ხალხისთვის"
$ws.Range("E76").Value = "საარჩევნო ბლოკი პაატა ბურჭულაძე სახელმწიფო ხალხისთვის"

# --- Row heights for the newly appended rows (auto-fit look of the source) -
$ws.Rows.Item(60).RowHeight = 40.25
$ws.Rows.Item(61).RowHeight = 40.25
$ws.Rows.Item(62).RowHeight = 40.25
$ws.Rows.Item(63).RowHeight = 14.7
$ws.Rows.Item(64).RowHeight = 14.7
$ws.Rows.Item(65).RowHeight = 14.7
$ws.Rows.Item(66).RowHeight = 14.7
$ws.Rows.Item(67).RowHeight = 14.7
$ws.Rows.Item(68).RowHeight = 27.45
$ws.Rows.Item(69).RowHeight = 12.8
$ws.Rows.Item(70).RowHeight = 12.8
$ws.Rows.Item(71).RowHeight = 12.8
$ws.Rows.Item(72).RowHeight = 12.8
$ws.Rows.Item(73).RowHeight = 12.8
$ws.Rows.Item(74).RowHeight = 12.8
$ws.Rows.Item(75).RowHeight = 12.8
$ws.Rows.Item(76).RowHeight = 12.8

# --- View state: scroll position + active selection moved with the data ---
$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("C71").Select()
